# Add 6 new "Summary_NN" sheets (Summary_58 .. Summary_63) at the end of the
# workbook, following the exact repeating template already used by the
# existing Summary_* sheets (header row + COUNTA/COUNTA/ratio formulas).

$wb = $excel.ActiveWorkbook

# Common header row used by every summary sheet.
$headers = @("Column Heading", "Count", "Total", "Percentage")

function Set-HeaderRow($ws) {
    $ws.Range("A1").Value = $headers[0]
    $ws.Range("B1").Value = $headers[1]
    $ws.Range("C1").Value = $headers[2]
    $ws.Range("D1").Value = $headers[3]
}

# Rows 2-6: one row per TestData column (Name/Age/City/Score/Comments),
# referencing TestData!A:A .. TestData!E:E for the "Count" column and always
# TestData!A:A for "Total".
function Set-TestDataRows($ws) {
    $labels = @("Name", "Age", "City", "Score", "Comments")
    $cols = @("A", "B", "C", "D", "E")
    for ($i = 0; $i -lt 5; $i++) {
        $r = $i + 2
        $ws.Range("A$r").Value = $labels[$i]
        $ws.Range("B$r").Formula = "=COUNTA(TestData!$($cols[$i]):$($cols[$i]))-1"
        $ws.Range("C$r").Formula = "=COUNTA(TestData!A:A)-1"
        $ws.Range("D$r").Formula = "=B$r/C$r"
        $ws.Range("D$r").NumberFormat = "0.00%"
    }
}

# Extra row 7 (no label in column A) referencing TestData!F:F.
function Set-ExtraTestDataRow7($ws) {
    $ws.Range("B7").Formula = "=COUNTA(TestData!F:F)-1"
    $ws.Range("C7").Formula = "=COUNTA(TestData!A:A)-1"
    $ws.Range("D7").Formula = "=B7/C7"
    $ws.Range("D7").NumberFormat = "0.00%"
}

# Rows 2-5: one row per Summary_1 column (Column Heading/Count/Total/
# Percentage), referencing Summary_1!A:A .. Summary_1!D:D, plus row 6 (no
# label) referencing Summary_1!E:E.
function Set-Summary1Rows($ws) {
    $labels = @("Column Heading", "Count", "Total", "Percentage")
    $cols = @("A", "B", "C", "D")
    for ($i = 0; $i -lt 4; $i++) {
        $r = $i + 2
        $ws.Range("A$r").Value = $labels[$i]
        $ws.Range("B$r").Formula = "=COUNTA(Summary_1!$($cols[$i]):$($cols[$i]))-1"
        $ws.Range("C$r").Formula = "=COUNTA(Summary_1!A:A)-1"
        $ws.Range("D$r").Formula = "=B$r/C$r"
        $ws.Range("D$r").NumberFormat = "0.00%"
    }
    $ws.Range("B6").Formula = "=COUNTA(Summary_1!E:E)-1"
    $ws.Range("C6").Formula = "=COUNTA(Summary_1!A:A)-1"
    $ws.Range("D6").Formula = "=B6/C6"
    $ws.Range("D6").NumberFormat = "0.00%"
}

function Add-SheetAtEnd($name) {
    $lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
    $ws = $wb.Worksheets.Add($null, $lastSheet)
    $ws.Name = $name
    return $ws
}

# sheet60 -> Summary_58 : header + TestData rows (dimension A1:D6)
$ws = Add-SheetAtEnd "Summary_58"
Set-HeaderRow $ws
Set-TestDataRows $ws

# sheet61 -> Summary_59 : header only (dimension A1:D1)
$ws = Add-SheetAtEnd "Summary_59"
Set-HeaderRow $ws

# sheet62 -> Summary_60 : header + TestData rows + extra col-F row (dimension A1:D7)
$ws = Add-SheetAtEnd "Summary_60"
Set-HeaderRow $ws
Set-TestDataRows $ws
Set-ExtraTestDataRow7 $ws

# sheet63 -> Summary_61 : header + Summary_1 rows (dimension A1:D6)
$ws = Add-SheetAtEnd "Summary_61"
Set-HeaderRow $ws
Set-Summary1Rows $ws

# sheet64 -> Summary_62 : header + TestData rows (dimension A1:D6)
$ws = Add-SheetAtEnd "Summary_62"
Set-HeaderRow $ws
Set-TestDataRows $ws

# sheet65 -> Summary_63 : header only (dimension A1:D1)
$ws = Add-SheetAtEnd "Summary_63"
Set-HeaderRow $ws

# Restore the original active sheet/selection (Summary_1) so the workbook
# view state matches the pre-edit file (adding sheets activates the new one).
$wb.Worksheets.Item("Summary_1").Activate()
$wb.Worksheets.Item("Summary_1").Range("A1").Select()
